$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 (CasesTab query) dropped its trailing "Cohort" output column -- the
# OPTIONAL MATCH (co:cohort) / WITH ... co plumbing stays, only the final
# RETURN line (and the now-unneeded trailing comma on the previous line)
# goes away.
$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Pulmonary Carcinoma"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesTabQuery.TrimEnd("`r", "`n")

# Author re-selected B2 (the cell they just edited) before saving, scrolling
# the view back to the top -- the previous selection had left the window
# scrolled down to B4.
$null = $ws.Range("B2").Select()
